$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 287-288, shifting the existing data (old rows
# 287-380) down to 289-382. This mirrors the diff where two brand-new
# "Choclero" records (dated 44559) are inserted and every subsequent row's
# content shifts down by two positions.
$ws.Rows("287:288").Insert()

# Populate the first new row (287): Choclero / Región Metropolitana record.
$ws.Cells.Item(287, 1).Value = 9
$ws.Cells.Item(287, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(287, 3).Value = "Metropolitana"
$ws.Cells.Item(287, 4).Value = 44559
$ws.Cells.Item(287, 5).Value = 13
$ws.Cells.Item(287, 6).Value = 100112024
$ws.Cells.Item(287, 7).Value = "Choclo"
$ws.Cells.Item(287, 8).Value = "Choclero"
$ws.Cells.Item(287, 9).Value = "Primera"
$ws.Cells.Item(287, 10).Value = 7900
$ws.Cells.Item(287, 11).Value = 300
$ws.Cells.Item(287, 12).Value = 350
$ws.Cells.Item(287, 13).Value = 325
$ws.Cells.Item(287, 14).Value = "$/unidad"
$ws.Cells.Item(287, 15).Value = "Región Metropolitana"
$ws.Cells.Item(287, 16).Value = 325
$ws.Cells.Item(287, 17).Value = 1
$ws.Cells.Item(287, 18).Value = "Hortaliza"

# Populate the second new row (288): Choclero / Región Metropolitana record.
$ws.Cells.Item(288, 1).Value = 9
$ws.Cells.Item(288, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(288, 3).Value = "Metropolitana"
$ws.Cells.Item(288, 4).Value = 44559
$ws.Cells.Item(288, 5).Value = 13
$ws.Cells.Item(288, 6).Value = 100112024
$ws.Cells.Item(288, 7).Value = "Choclo"
$ws.Cells.Item(288, 8).Value = "Choclero"
$ws.Cells.Item(288, 9).Value = "Primera"
$ws.Cells.Item(288, 10).Value = 5200
$ws.Cells.Item(288, 11).Value = 200
$ws.Cells.Item(288, 12).Value = 250
$ws.Cells.Item(288, 13).Value = 225
$ws.Cells.Item(288, 14).Value = "$/unidad"
$ws.Cells.Item(288, 15).Value = "Región Metropolitana"
$ws.Cells.Item(288, 16).Value = 225
$ws.Cells.Item(288, 17).Value = 1
$ws.Cells.Item(288, 18).Value = "Hortaliza"
